$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 20 ("RandomForestGroupVarImp.R") - the sheet is shrinking
# by one data row (A1:B25 -> A1:B24), and everything below shifts up.
$ws.Rows.Item(20).Delete()

# Row 18 becomes the new "pair-wise" sources-of-variance entry, with taller row height.
$ws.Range("A18").Value2 = "GenFigSourcesOfVarPairWise.R"
$ws.Range("B18").Value2 = "Manually checked that the number of technical replicate pairs per lab was correct. The results in all looks as expected. No further tests were deemed necessary."
$ws.Rows.Item(18).RowHeight = 30

# Row 19's file/function is now the (non-fine-grained) group var-imp script.
$ws.Range("A19").Value2 = "RandomForestGroupVarImp.R"

# Restore the selection to reflect the edited area.
$ws.Range("A20").Select()
